# fix(CWL): use `SafeQueryTypes` for decltype
# Adds a new "cwl_warn_decltype_missing" row (row 53) to the General sheet,
# mirroring the layout/styles of the other cwl_warn_* rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# --- Create row 53 with the same row-level styling as row 52 (the last data row) ---
$ws.Rows(52).Copy()
$ws.Rows(53).Insert(-4121)          # xlShiftDown

# --- Fix up per-cell styles to match the other "cwl_warn_*" rows (e.g. row 51) ---
# Column A uses style "2" (bold/green key style) on cwl_warn_* rows, column B is blank
# with style "9" on those same rows.
$ws.Range("A51").Copy()
$ws.Range("A53").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B51").Copy()
$ws.Range("B53").PasteSpecial(-4122) # xlPasteFormats

# --- Values ---
$ws.Range("A53").Value = "cwl_warn_decltype_missing"
$ws.Range("B53").Value = ""

$ws.Range("C53").Value = "failed to query decltype from mod: {0} -> {1}" + $nl + "it might've failed to load or is missing (transitive) dependencies" + $nl + "this is not an exception from CWL"

$zh = "无法查询MOD声明类型" + ": {0} -> {1}  " + $nl + "或许它未能加载或缺少（传递）依赖项" + "  " + $nl + "这并不是一个 " + "CWL " + "异常"
$ws.Range("D53").Value = $zh

# --- Apply the mixed-font rich text formatting to D53, matching the CJK / code-font
#     split used throughout the rest of the sheet (微软雅黑 for CJK prose, Cascadia
#     Code for the code-ish / latin segments). ---
$d = $ws.Range("D53")

$r1 = $d.Characters(1, 11)   # 无法查询MOD声明类型
$r1.Font.ColorIndex = -4105
$r1.Font.Name = "微软雅黑"

$r2 = $d.Characters(12, 15)  # : {0} -> {1}  \n
$r2.Font.ColorIndex = -4105
$r2.Font.Name = "Cascadia Code"

$r3 = $d.Characters(27, 17)  # 或许它未能加载或缺少（传递）依赖项
$r3.Font.ColorIndex = -4105
$r3.Font.Name = "微软雅黑"

$r4 = $d.Characters(44, 3)   #   \n
$r4.Font.ColorIndex = -4105
$r4.Font.Name = "Cascadia Code"

$r5 = $d.Characters(47, 7)   # 这并不是一个 
$r5.Font.ColorIndex = -4105
$r5.Font.Name = "微软雅黑"

$r6 = $d.Characters(54, 4)   # CWL 
$r6.Font.ColorIndex = -4105
$r6.Font.Name = "Cascadia Code"

$r7 = $d.Characters(58, 2)   # 异常
$r7.Font.ColorIndex = -4105
$r7.Font.Name = "微软雅黑"

# --- Row height: the new row wraps to roughly the same height as the other long,
#     multi-line cwl_warn_* rows. ---
$ws.Rows(53).RowHeight = 93

# --- Update the view: scroll back to the top-left of the data area and leave the
#     selection where Excel would land after entering this row. ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 46
$null = $ws.Range("D56").Select()
